# Update panel2 map for auto parents, and add a new "panel2_v2" sheet
# that captures the auto->manual parent chain (only the rows that
# actually participate in the chain).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("panel2")

# --- panel2: backfill the "RealAutoParent" (column B) for the rows that
# feed the auto-gating chain. They were all "NA" before. ---
$ws2.Range("B5").Value = "root"

$ws2.Range("C10").Copy()
$ws2.Range("B10").PasteSpecial(-4122)
$ws2.Range("B10").Value = "Live immune cells (CD45+ PE-)"

$ws2.Range("C12").Copy()
$ws2.Range("B12").PasteSpecial(-4122)
$ws2.Range("B12").Value = "Live Single immune cells(FSC-H/FSC-W)"

# panel2 is no longer the selected/active tab; move the cursor first.
$ws2.Activate()
$ws2.Range("C12").Select()

# --- add the new "panel2_v2" sheet, after panel2, with the condensed
# Auto/RealAutoParent/Manual mapping. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "panel2_v2"

$ws3.Range("A1").Value = "Auto"
$ws3.Range("B1").Value = "RealAutoParent"
$ws3.Range("C1").Value = "Manual"

$ws3.Range("A2").Value = "PE-A"
$ws3.Range("B2").Value = "root"
$ws3.Range("C2").Value = "Live immune cells (CD45+ PE-)"

$ws3.Range("A3").Value = "Singlets"
$ws3.Range("B3").Value = "Live immune cells (CD45+ PE-)"
$ws3.Range("C3").Value = "Live Single immune cells(FSC-H/FSC-W)"

$ws3.Range("A4").Value = "PBMC"
$ws3.Range("B4").Value = "Live Single immune cells(FSC-H/FSC-W)"
$ws3.Range("C4").Value = "Live Single PBMCs (SSC-A/FSC-A)"

$ws3.Range("A5").Value = "D_NK_M"
$ws3.Range("B5").Value = "Live Single PBMCs (SSC-A/FSC-A)"
$ws3.Range("C5").Value = "DC NK MONOCYTES (CD3- CD19-)"

$ws3.Range("A6").Value = "CD14+"
$ws3.Range("B6").Value = "DC NK MONOCYTES (CD3- CD19-)"
$ws3.Range("C6").Value = "MONOCYTES (CD14+)"

$ws3.Range("A7").Value = "CD14+/CD16+"
$ws3.Range("B7").Value = "MONOCYTES (CD14+)"
$ws3.Range("C7").Value = "Non classical monocytes (CD16+ CD14+)"

$ws3.Range("A8").Value = "CD16-"
$ws3.Range("B8").Value = "MONOCYTES (CD14+)"
$ws3.Range("C8").Value = "Classical monocytes (CD16- CD14+)"

$ws3.Range("A9").Value = "CD20-"
$ws3.Range("B9").Value = "DC NK MONOCYTES (CD3- CD19-)"
$ws3.Range("C9").Value = "DC NK (CD20- CD14-)"

$ws3.Range("A10").Value = "Dendritic"
$ws3.Range("B10").Value = "DC NK (CD20- CD14-)"
$ws3.Range("C10").Value = "DC (HLA-DR+)"

$ws3.Range("A11").Value = "BB515-A+BV 711-A-"
$ws3.Range("B11").Value = "DC (HLA-DR+)"
$ws3.Range("C11").Value = "Myeloid DC (CD11c+ CD123-)"

$ws3.Range("A12").Value = "BB515-A-BV 711-A-"
$ws3.Range("B12").Value = "DC (HLA-DR+)"
$ws3.Range("C12").Value = "Plasmacytoid DC (CD11c- CD123+)"

$ws3.Range("A13").Value = "CD20-/CD16+"
$ws3.Range("B13").Value = "DC NK (CD20- CD14-)"
$ws3.Range("C13").Value = "NK (CD16+)"

$ws3.Range("A14").Value = "CD16+CD56+"
$ws3.Range("B14").Value = "NK (CD16+)"
$ws3.Range("C14").Value = "NK CD56LO"

$ws3.Range("A15").Value = "CD56PlusPlus"
$ws3.Range("B15").Value = "NK (CD16+)"
$ws3.Range("C15").Value = "NK CD56HI"

# panel2_v2 becomes the active tab, with its own selection.
$ws3.Activate()
$ws3.Range("H26").Select()
